$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update revised figures for rows 218-222 (columns B,C,E,F,G,H; D unchanged) ---
$ws.Range("B218").Value = 56697
$ws.Range("C218").Value = 4805
$ws.Range("E218").Value = 4735
$ws.Range("F218").Value = 51892
$ws.Range("G218").Value = 3927
$ws.Range("H218").Value = 47965

$ws.Range("B219").Value = 57221
$ws.Range("C219").Value = 4910
$ws.Range("E219").Value = 4840
$ws.Range("F219").Value = 52311
$ws.Range("G219").Value = 4068
$ws.Range("H219").Value = 48244

$ws.Range("B220").Value = 57157
$ws.Range("C220").Value = 4749
$ws.Range("E220").Value = 4679
$ws.Range("F220").Value = 52408
$ws.Range("G220").Value = 4165
$ws.Range("H220").Value = 48243

$ws.Range("B221").Value = 58214
$ws.Range("C221").Value = 4846
$ws.Range("E221").Value = 4775
$ws.Range("F221").Value = 53369
$ws.Range("G221").Value = 4170
$ws.Range("H221").Value = 49198

$ws.Range("B222").Value = 57870
$ws.Range("C222").Value = 4828
$ws.Range("E222").Value = 4757
$ws.Range("F222").Value = 53042
$ws.Range("G222").Value = 3932
$ws.Range("H222").Value = 49110

# --- Append the new month row (223) : 01-06-2021 ---
# Column A must end up as a plain text shared string (like the other period
# labels in the sheet), not an auto-recognized date serial. We build the
# label via a formula (so the "looks like a date" literal-input heuristic
# never sees the raw text), force-calculate it, then convert that formula
# to its static text result with a values-only paste.
$ws.Range("A223").Formula = '=TRIM("01-06-2021 ")'
$excel.CalculateFull()
$ws.Range("A223").Copy()
$ws.Range("A223").PasteSpecial(-4163)

$ws.Range("B223").Value = 56685
$ws.Range("C223").Value = 4769
$ws.Range("D223").Value = 71
$ws.Range("E223").Value = 4698
$ws.Range("F223").Value = 51916
$ws.Range("G223").Value = 3497
$ws.Range("H223").Value = 48419
